# Apply cryptos list update (prices + volume%) from Aug 13 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.780.04"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "2.626.87"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'520.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").Value = "'144.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").Value = "2.635.72"
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "'0.126"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "3.087.01"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "58.758.79"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "'20.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "2.642.33"
$ws.Range("E18").Value = "  +2.45%  "
$ws.Range("D19").Value = "'344.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "'4.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "'10.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("D22").Value = "'6.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'61.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "'0.414"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").Value = "'0.163"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.28%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "0.0₃0796"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").Value = "'7.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'6.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'18.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "'150.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").Value = "'0.977"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.46%  "
$ws.Range("D36").Value = "'3.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "'36.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("D39").Value = "'0.835"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("E40").Value = "  +1.64%  "
$ws.Range("D41").Value = "'1.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("D42").Value = "'277.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.24%  "
$ws.Range("D43").Value = "'0.996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "'0.0981"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.599"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'19.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'10.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0520"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.93%  "
$ws.Range("D49").Value = "1.989.86"
$ws.Range("E49").Value = "  +3.34%  "
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "'4.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.21%  "
